$wb = $excel.ActiveWorkbook

# This script applies refreshed market-price values (columns H-N) across all
# eight Leve-profit worksheets, matching the upstream scheduled-runner update.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1009.55
$ws.Range("I28").Value = 957.4211
$ws.Range("K28").Value = 957.4211
$ws.Range("M28").Value = -472.4211
$ws.Range("H40").Value = 4268.0625
$ws.Range("J40").Value = 4592.0713
$ws.Range("L40").Value = 4592.0713
$ws.Range("N40").Value = -4942.0713
$ws.Range("H43").Value = 4293.9414
$ws.Range("I43").Value = 4164
$ws.Range("J43").Value = 4384.9
$ws.Range("K43").Value = 4164
$ws.Range("L43").Value = 4384.9
$ws.Range("M43").Value = -4095
$ws.Range("N43").Value = -4522.9
$ws.Range("H99").Value = 1221.4615
$ws.Range("I99").Value = 1486.7778
$ws.Range("J99").Value = 624.5
$ws.Range("K99").Value = 4460.3334
$ws.Range("L99").Value = 1873.5
$ws.Range("M99").Value = -2962.3334
$ws.Range("N99").Value = -4869.5
$ws.Range("H115").Value = 1000.8182
$ws.Range("I115").Value = 1000.8182
$ws.Range("K115").Value = 3002.4546
$ws.Range("M115").Value = -1435.4546
$ws.Range("H116").Value = 3134.5
$ws.Range("I116").Value = 2917
$ws.Range("K116").Value = 2917
$ws.Range("M116").Value = 525
$ws.Range("H132").Value = 34490700
$ws.Range("I132").Value = 37041120
$ws.Range("K132").Value = 111123360
$ws.Range("M132").Value = -111120830
$ws.Range("H135").Value = 5477
$ws.Range("I135").Value = 2667.7144
$ws.Range("J135").Value = 15309.5
$ws.Range("K135").Value = 24009.4296
$ws.Range("L135").Value = 137785.5
$ws.Range("M135").Value = -21474.4296
$ws.Range("N135").Value = -142855.5
$ws.Range("H138").Value = 314424.03
$ws.Range("I138").Value = 760.7222
$ws.Range("J138").Value = 717705.4399999999
$ws.Range("K138").Value = 2282.1666
$ws.Range("L138").Value = 2153116.32
$ws.Range("M138").Value = 2857.8334
$ws.Range("N138").Value = -2163396.32

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 978.41174
$ws.Range("J2").Value = 1811
$ws.Range("L2").Value = 1811
$ws.Range("N2").Value = -2037
$ws.Range("H116").Value = 978.41174
$ws.Range("J116").Value = 1811
$ws.Range("L116").Value = 1811
$ws.Range("N116").Value = -6399
$ws.Range("H122").Value = 3930.4634
$ws.Range("I122").Value = 3468.125
$ws.Range("K122").Value = 10404.375
$ws.Range("M122").Value = -7954.375

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 978.41174
$ws.Range("J3").Value = 1811
$ws.Range("L3").Value = 1811
$ws.Range("N3").Value = -2039
$ws.Range("H38").Value = 12218.923
$ws.Range("I38").Value = 7225.6
$ws.Range("J38").Value = 28863.334
$ws.Range("K38").Value = 7225.6
$ws.Range("L38").Value = 28863.334
$ws.Range("M38").Value = -6809.6
$ws.Range("N38").Value = -29695.334
$ws.Range("H94").Value = 705.9231
$ws.Range("I94").Value = 437.85715
$ws.Range("J94").Value = 1831.8
$ws.Range("K94").Value = 437.85715
$ws.Range("L94").Value = 1831.8
$ws.Range("M94").Value = 13.14285000000001
$ws.Range("N94").Value = -2733.8
$ws.Range("H107").Value = 2340.3
$ws.Range("I107").Value = 1850.7
$ws.Range("J107").Value = 3319.5
$ws.Range("K107").Value = 1850.7
$ws.Range("L107").Value = 3319.5
$ws.Range("M107").Value = 69.29999999999995
$ws.Range("N107").Value = -7159.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 416.33334
$ws.Range("I7").Value = 196
$ws.Range("J7").Value = 592.6
$ws.Range("K7").Value = 196
$ws.Range("L7").Value = 592.6
$ws.Range("M7").Value = -83
$ws.Range("N7").Value = -818.6
$ws.Range("H22").Value = 305.1111
$ws.Range("I22").Value = 295.2
$ws.Range("J22").Value = 354.66666
$ws.Range("K22").Value = 295.2
$ws.Range("L22").Value = 354.66666
$ws.Range("M22").Value = 54.80000000000001
$ws.Range("N22").Value = -1054.66666
$ws.Range("H45").Value = 25000
$ws.Range("J45").Value = 25000
$ws.Range("L45").Value = 25000
$ws.Range("N45").Value = -26186
$ws.Range("H107").Value = 1718.0476
$ws.Range("J107").Value = 3424
$ws.Range("L107").Value = 3424
$ws.Range("N107").Value = -7264
$ws.Range("H132").Value = 1113113.8
$ws.Range("I132").Value = 1381048.5
$ws.Range("K132").Value = 4143145.5
$ws.Range("M132").Value = -4140615.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9139356
$ws.Range("I4").Value = 11425123
$ws.Range("K4").Value = 34275369
$ws.Range("M4").Value = -34275257
$ws.Range("H5").Value = 5685.25
$ws.Range("J5").Value = 8305.333000000001
$ws.Range("L5").Value = 24915.999
$ws.Range("N5").Value = -25139.999
$ws.Range("H12").Value = 507.125
$ws.Range("J12").Value = 463.8
$ws.Range("L12").Value = 1391.4
$ws.Range("N12").Value = -1737.4
$ws.Range("H135").Value = 5685.25
$ws.Range("J135").Value = 8305.333000000001
$ws.Range("L135").Value = 74747.997
$ws.Range("N135").Value = -79817.997

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 29998
$ws.Range("J26").Value = 29998
$ws.Range("L26").Value = 29998
$ws.Range("N26").Value = -30558
$ws.Range("H31").Value = 2635
$ws.Range("J31").Value = 2635
$ws.Range("L31").Value = 2635
$ws.Range("N31").Value = -3219
$ws.Range("H37").Value = 2635
$ws.Range("J37").Value = 2635
$ws.Range("L37").Value = 2635
$ws.Range("N37").Value = -3189
$ws.Range("H48").Value = 28998
$ws.Range("J48").Value = 28998
$ws.Range("L48").Value = 28998
$ws.Range("N48").Value = -29968
$ws.Range("H50").Value = 29998
$ws.Range("J50").Value = 29998
$ws.Range("L50").Value = 29998
$ws.Range("N50").Value = -30994
$ws.Range("H107").Value = 591.46155
$ws.Range("I107").Value = 562.8
$ws.Range("K107").Value = 562.8
$ws.Range("M107").Value = 1357.2
$ws.Range("H113").Value = 3344.889
$ws.Range("I113").Value = 2019.6
$ws.Range("K113").Value = 2019.6
$ws.Range("M113").Value = 150.4000000000001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4030
$ws.Range("J16").Value = 4030
$ws.Range("L16").Value = 4030
$ws.Range("N16").Value = -4370
$ws.Range("H22").Value = 2949.4167
$ws.Range("I22").Value = 900
$ws.Range("K22").Value = 900
$ws.Range("M22").Value = -605
$ws.Range("H27").Value = 2949.4167
$ws.Range("I27").Value = 900
$ws.Range("K27").Value = 900
$ws.Range("M27").Value = -793
$ws.Range("H40").Value = 4034.4546
$ws.Range("I40").Value = 4544.263
$ws.Range("J40").Value = 3342.5715
$ws.Range("K40").Value = 4544.263
$ws.Range("L40").Value = 3342.5715
$ws.Range("M40").Value = -4408.263
$ws.Range("N40").Value = -3614.5715
$ws.Range("H41").Value = 30000
$ws.Range("J41").Value = 30000
$ws.Range("L41").Value = 30000
$ws.Range("N41").Value = -30876
$ws.Range("H55").Value = 450.86667
$ws.Range("I55").Value = 385.33334
$ws.Range("J55").Value = 549.1667
$ws.Range("K55").Value = 385.33334
$ws.Range("L55").Value = 549.1667
$ws.Range("M55").Value = -212.33334
$ws.Range("N55").Value = -895.1667
$ws.Range("H61").Value = 24636.908
$ws.Range("I61").Value = 27369.21
$ws.Range("K61").Value = 27369.21
$ws.Range("M61").Value = -27167.21
$ws.Range("H101").Value = 49738
$ws.Range("J101").Value = 49738
$ws.Range("L101").Value = 49738
$ws.Range("N101").Value = -56228
$ws.Range("H113").Value = 24636.908
$ws.Range("I113").Value = 27369.21
$ws.Range("K113").Value = 27369.21
$ws.Range("M113").Value = -25199.21

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 39965
$ws.Range("J50").Value = 39965
$ws.Range("L50").Value = 39965
$ws.Range("N50").Value = -41227
